# Generate Report for Handoff
# Adds two new localized files (80e49444-... and af18a52e-...) as rows 4 & 5
# across the Overview / zh-cn / de-de sheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$mdCommit   = "f4d06d6e9d1b7b69d2d87778346ab7134c387149"
$zhcnCommit = "74364f1932ee4149d6af2bf0eef6f2c3e7467d09"
$dedeCommit = "9ad784043caf788e9016ad3b1fe5f22719b38abc"

function MdUrl($guid) {
    return "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$guid.md"
}
function ZhCnXlfUrl($fname) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhcnCommit/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$fname"
}
function DeDeXlfUrl($fname) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$dedeCommit/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$fname"
}

$guid1 = "80e49444-9e42-4422-bfc9-ffd54ea8cdff"
$hash1 = "74521c2fd6e93a574f80ce8fba105456cd4d8a56"
$guid2 = "af18a52e-2a81-4c18-8253-40ba07ba2e92"
$hash2 = "60c9dffd4504ebbb57fd6eee2ababab60620bcc8"

$zhcnXlf1 = "$guid1.$hash1.zh-cn.xlf"
$zhcnXlf2 = "$guid2.$hash2.zh-cn.xlf"
$dedeXlf1 = "$guid1.$hash1.de-de.xlf"
$dedeXlf2 = "$guid2.$hash2.de-de.xlf"

$status = "Ready for handoff"
$zhcnHandoffDatetime = "2016-03-31 05:11:08"
$dedeHandoffDatetime = "2016-03-31 05:11:19"
$latestHandoffDate   = "2016-03-31 05:11:19"
$emptyDatetime = "0001-01-01 00:00:00"
$mdExt = ".md"
$include = "Include"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "$guid1.md"
$wsOverview.Range("B4").Value = $status
$wsOverview.Range("C4").Value = $status
$wsOverview.Range("D4").Value = $latestHandoffDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), (MdUrl $guid1), "", "", "$guid1.md")

$wsOverview.Range("A5").Value = "$guid2.md"
$wsOverview.Range("B5").Value = $status
$wsOverview.Range("C5").Value = $status
$wsOverview.Range("D5").Value = $latestHandoffDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), (MdUrl $guid2), "", "", "$guid2.md")

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "$guid1.md"
$wsZhCn.Range("B4").Value = $mdExt
$wsZhCn.Range("C4").Value = $status
$wsZhCn.Range("D4").Value = $zhcnXlf1
$wsZhCn.Range("E4").Value = $zhcnHandoffDatetime
$wsZhCn.Range("H4").Value = $emptyDatetime
$wsZhCn.Range("J4").Value = $include
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), (MdUrl $guid1), "", "", "$guid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D4"), (ZhCnXlfUrl $zhcnXlf1), "", "", $zhcnXlf1)

$wsZhCn.Range("A5").Value = "$guid2.md"
$wsZhCn.Range("B5").Value = $mdExt
$wsZhCn.Range("C5").Value = $status
$wsZhCn.Range("D5").Value = $zhcnXlf2
$wsZhCn.Range("E5").Value = $zhcnHandoffDatetime
$wsZhCn.Range("H5").Value = $emptyDatetime
$wsZhCn.Range("J5").Value = $include
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), (MdUrl $guid2), "", "", "$guid2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D5"), (ZhCnXlfUrl $zhcnXlf2), "", "", $zhcnXlf2)

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "$guid1.md"
$wsDeDe.Range("B4").Value = $mdExt
$wsDeDe.Range("C4").Value = $status
$wsDeDe.Range("D4").Value = $dedeXlf1
$wsDeDe.Range("E4").Value = $latestHandoffDate
$wsDeDe.Range("H4").Value = $emptyDatetime
$wsDeDe.Range("J4").Value = $include
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), (MdUrl $guid1), "", "", "$guid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D4"), (DeDeXlfUrl $dedeXlf1), "", "", $dedeXlf1)

$wsDeDe.Range("A5").Value = "$guid2.md"
$wsDeDe.Range("B5").Value = $mdExt
$wsDeDe.Range("C5").Value = $status
$wsDeDe.Range("D5").Value = $dedeXlf2
$wsDeDe.Range("E5").Value = $latestHandoffDate
$wsDeDe.Range("H5").Value = $emptyDatetime
$wsDeDe.Range("J5").Value = $include
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), (MdUrl $guid2), "", "", "$guid2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D5"), (DeDeXlfUrl $dedeXlf2), "", "", $dedeXlf2)
